$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the three features that shipped in 1.5.0 as completed in that version.
# (Column C = "Completed Version")
$ws.Range("C4").Value = "1.5.0"
$ws.Range("C6").Value = "1.5.0"
$ws.Range("C7").Value = "1.5.0"

# Add the newly-completed "Restore history button" feature as row 12.
$ws.Range("A12").Value = "Restore history button"
$ws.Range("B12").Value = "Removing all the history on accident can be a pain. Is it possible to restore that with a button somewhere?"
$ws.Range("C12").Value = "1.5.1"
$ws.Range("D12").Value = "Bart van den Hoek - Store Review"

# Filter the tracker down to the still-outstanding (blank "Completed Version") rows.
$ws.Range("A1:E12").AutoFilter(3, @(""), 7)

# Register the hidden _FilterDatabase defined name Excel creates for the filter.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$1:`$E`$12")
$filterName.Visible = $false

# Leave the selection where the author last left it.
$ws.Range("B17").Select()
